# Apply "repull data, push all data, mean calculation" changes:
# Update dSF (column F) values for a handful of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -6
$ws.Range("F8").Value = -8
$ws.Range("F11").Value = -5
$ws.Range("F13").Value = -5
